# date format changed in excel sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New dd/mm/yyyy text values for the dateSanctioned column (G2:G20)
$dates = @(
    "11/03/2013",
    "14/05/2013",
    "17/07/2013",
    "19/09/2013",
    "22/11/2013",
    "25/01/2014",
    "30/03/2014",
    "02/06/2014",
    "05/08/2014",
    "08/10/2014",
    "11/12/2014",
    "13/02/2015",
    "18/04/2015",
    "21/06/2015",
    "24/08/2015",
    "27/10/2015",
    "30/12/2015",
    "03/03/2016",
    "06/05/2016"
)

# Switch the whole column G to a text number format (was a date format)
$colG = $ws.Range("G1:G20")
$colG.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $dates[$i]
}

# Remove the trailing blank rows (21-23) that only held leftover formatting
$ws.Rows.Item(21).Delete() | Out-Null
$ws.Rows.Item(21).Delete() | Out-Null
$ws.Rows.Item(21).Delete() | Out-Null

# Reset the active selection to I1
$ws.Range("I1").Select() | Out-Null
